$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID,DRUG,ON_RECEPT,NO_PACKAGES_AVAILABLE,DATE,RECEPT_ID"
$ws.Range("A2").Value = "1,POLOPIRYNA,NO,1200,2024-03-25 00:00:00,nan"
$ws.Range("A3").Value = "3,AMOTAX,YES,336,2015-02-14 00:00:00,nan"
$ws.Range("A4").Value = "5,GSGDF,YES,455,2025-05-24 19:55:06,nan"
$ws.Range("A5").Value = "6,JHV,YES,4445,2025-05-25 09:19:09,nan"
$ws.Range("A6").Value = "7,RREWRW,NO,4432,2025-05-25 09:19:19,nan"
$ws.Range("A7").Value = "10,EWQ,YES,232,2025-05-26 21:05:04,nan"
$ws.Range("A8").Value = "12,SDSAD,YES,2311,2025-05-27 15:51:10,2344.0"
